$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="61×89=5429"; New="93×54=5022"},
    @{Row=1;  Col=2; Old="33×70=2310"; New="94×43=4042"},
    @{Row=1;  Col=3; Old="78×71=5538"; New="44×58=2552"},
    @{Row=1;  Col=4; Old="11×56=616";  New="12×55=660"},
    @{Row=1;  Col=5; Old="90×15=1350"; New="55×36=1980"},

    @{Row=5;  Col=1; Old="62×47=2914"; New="40×38=1520"},
    @{Row=5;  Col=2; Old="25×11=275";  New="18×43=774"},
    @{Row=5;  Col=3; Old="30×59=1770"; New="56×36=2016"},
    @{Row=5;  Col=4; Old="37×34=1258"; New="77×96=7392"},
    @{Row=5;  Col=5; Old="12×70=840";  New="40×27=1080"},

    @{Row=10; Col=1; Old="15×97=1455"; New="95×32=3040"},
    @{Row=10; Col=2; Old="77×91=7007"; New="89×44=3916"},
    @{Row=10; Col=3; Old="45×82=3690"; New="49×46=2254"},
    @{Row=10; Col=4; Old="27×70=1890"; New="26×37=962"},
    @{Row=10; Col=5; Old="29×75=2175"; New="69×65=4485"},

    @{Row=15; Col=1; Old="78×84=6552"; New="84×44=3696"},
    @{Row=15; Col=2; Old="14×36=504";  New="51×64=3264"},
    @{Row=15; Col=3; Old="73×76=5548"; New="59×98=5782"},
    @{Row=15; Col=4; Old="96×71=6816"; New="92×57=5244"},
    @{Row=15; Col=5; Old="77×91=7007"; New="20×70=1400"},

    @{Row=20; Col=1; Old="21×90=1890"; New="65×93=6045"},
    @{Row=20; Col=2; Old="45×50=2250"; New="91×66=6006"},
    @{Row=20; Col=3; Old="35×42=1470"; New="50×28=1400"},
    @{Row=20; Col=4; Old="72×57=4104"; New="16×30=480"},
    @{Row=20; Col=5; Old="71×87=6177"; New="58×86=4988"}
)

foreach ($r in $replacements) {
    $cellRange = $t.Cell($r.Row, $r.Col).Range
    $cellRange.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 1)
}
